$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style from B1 to C1, then set header text
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$ws.Range("C1").Value = "Ax_"

# Fill column C (rows 2-131) with Ax_ values
$values = @(
    0.02183846773339848, 0.02270708734910754, 0.02361898583241229, 0.02457629390816973, 0.02558124406583362,
    0.02663617503552584, 0.02774353641314986, 0.02890589343322503, 0.03012593188700063, 0.03140646318216045,
    0.0327504295390646, 0.0341609093168419, 0.03564112246082207, 0.03719443606072736, 0.03882437000666317,
    0.04053460272736675, 0.04232897699213733, 0.04421150575457011, 0.04618637801248482, 0.0482579646542689,
    0.05043082425724899, 0.05270970879854797, 0.05509956923321591, 0.05760556088814617, 0.06023304861331924,
    0.0629876116247734, 0.06587504796388388, 0.06890137849155775, 0.07207285032207433, 0.07539593959417823,
    0.07887735346338687, 0.08252403118791796, 0.08634314488640094, 0.09034209477851418, 0.09452851383970862,
    0.0989102565158152, 0.1034953964640188, 0.1082922180001886, 0.1133092060498886, 0.1185550336302262,
    0.1240385465911286, 0.1297687453276781, 0.1357547631586213, 0.1420058410511875, 0.1485312983589038,
    0.1553404992284082, 0.1624428143234935, 0.1698475775113208, 0.1775640371574115, 0.1856013016840466,
    0.1939682790650441, 0.2026736099329353, 0.2117255940893398, 0.2211321100897944, 0.2309005278552871,
    0.2410376141570991, 0.2515494310090365, 0.2624412270915558, 0.2737173224681603, 0.2853809870123595,
    0.2974343131445313, 0.3098780836821342, 0.3227116358351644, 0.3359327226292517, 0.3495373733106375,
    0.3635197545757143, 0.3778720347697695, 0.3925842535072673, 0.4076441994722695, 0.4230372994515147,
    0.4387465219228787, 0.454752298753054, 0.4710324687358843, 0.4875622468073545, 0.5043142227875714,
    0.5212583934040509, 0.538362231125875, 0.5555907929675706, 0.5729068718905337, 0.5902711927289979,
    0.6076426536922175, 0.6249786134500703, 0.6422352226071124, 0.6593677970349821, 0.6763312290997248,
    0.6930804313355352, 0.7095708056367056, 0.725758729635308, 0.7416020506721074, 0.7570605767369721,
    0.7720965530230101, 0.7866751123789689, 0.8007646880128469, 0.8143373773350518, 0.827369246846828,
    0.8398405694662308, 0.8517359875969845, 0.8630445975125538, 0.8737599531487861, 0.8838799900519639,
    0.8934068728791764, 0.9023467723543033, 0.9107095798114119, 0.9185085692903686, 0.9257600184956263,
    0.9324828007298374, 0.9386979601485927, 0.9444282823679473, 0.9496978716448743, 0.9545317446237099,
    0.958955449101193, 0.9629947145221771, 0.9666751390920042, 0.9700219165866659, 0.973059604249833,
    0.9758119316584863, 0.9783016491651689, 0.9805504135116315, 0.9825787074617135, 0.9844057898093281,
    0.9860496718562833, 0.9875271163905259, 0.9888536552920273, 0.9900436221109798, 0.9911101962660576,
    0.9920654558659634, 0.9929204365380463, 0.9936851940313892, 0.9943688687318274, 0.9949797505712356
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 2, 3).Value = $values[$i]
}

